$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "correlation 5 and 6 - commons-lang": refresh the DLOC/Metrics-5 data table
# with the new commons-lang correlation numbers (6 data rows instead of the
# previous 11), and drop the now-unused trailing rows.
# ---------------------------------------------------------------------------

# New correlation data for rows 4-9 (Versions / DLOC columns).
# A (Sl.No) and row 7/8's B/C values are unchanged from before.
$ws.Range("B4").Value = 2.6
$ws.Range("C4").Value = 3046.0

$ws.Range("B5").Value = 3.2
$ws.Range("C5").Value = 60259.0

$ws.Range("B6").Value = 3.3
$ws.Range("C6").Value = 2687.0

$ws.Range("C7").Value = 4368.0

$ws.Range("B9").Value = 3.1
$ws.Range("C9").Value = 36569.0
$ws.Range("B9").NumberFormat = "0.00"

# Metrics 5 (D column) is a shared "78+(0.01*C#)" formula; the live data only
# spans rows 4:9 now, so re-enter it over that shrunk range.
$ws.Range("D4:D9").Formula = "=78+(0.01*C4)"

# Rows 10-14 no longer hold any data (old versions 3.7-3.1 incl. the stray
# "=3.1" formula cell in B14) - clear their contents, keep the formatting.
$ws.Range("A10:D14").ClearContents()

# Visual refresh: the kept header + untouched legacy cells pick up
# Arial / automatic (theme) text color.
$ws.Range("B3:G3").Font.Name = "Arial"
$ws.Range("B3:G3").Font.ThemeColor = 1

$ws.Range("A4:A9").Font.Name = "Arial"
$ws.Range("A4:A9").Font.ThemeColor = 1

$ws.Range("B7:B8").Font.Name = "Arial"
$ws.Range("B7:B8").Font.ThemeColor = 1

$ws.Range("C8").Font.Name = "Arial"
$ws.Range("C8").Font.ThemeColor = 1
